# This script applies a weekly update to the "Fruta, Terminal Hortofrutícola
# Agro Chillán - Pera" dataset: it inserts two new rows of data for the most
# recent reporting date right after the fixed header block (row 14), pushing
# all the previously existing data rows (old rows 15-136) down by two rows
# (to new rows 17-138), and fills the two freshly inserted rows (15 and 16)
# with the new week's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 15. Excel will push the
# existing rows 15-136 down to 17-138 and carry over cell formatting
# (e.g. the date number format in column D) from the surrounding rows.
$ws.Rows("15:16").Insert()

# New row 15: "Especial" quality entry for the new reporting date.
$ws.Cells.Item(15, 1).Value = 7
$ws.Cells.Item(15, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(15, 3).Value = "Ñuble"
$ws.Cells.Item(15, 4).Value = 44473
$ws.Cells.Item(15, 5).Value = 16
$ws.Cells.Item(15, 6).Value = "Fruta"
$ws.Cells.Item(15, 7).Value = 100104
$ws.Cells.Item(15, 8).Value = "Frutos de pepita"
$ws.Cells.Item(15, 9).Value = 100104005
$ws.Cells.Item(15, 10).Value = "Pera"
$ws.Cells.Item(15, 11).Value = "Packham's Triumph"
$ws.Cells.Item(15, 12).Value = "Especial"
$ws.Cells.Item(15, 13).Value = 30
$ws.Cells.Item(15, 14).Value = 11000
$ws.Cells.Item(15, 15).Value = 11000
$ws.Cells.Item(15, 16).Value = 11000
$ws.Cells.Item(15, 17).Value = "`$/caja 16 kilos empedrada"
$ws.Cells.Item(15, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(15, 19).Value = 688
$ws.Cells.Item(15, 20).Value = 16

# New row 16: "Primera" quality entry for the same new reporting date.
$ws.Cells.Item(16, 1).Value = 7
$ws.Cells.Item(16, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(16, 3).Value = "Ñuble"
$ws.Cells.Item(16, 4).Value = 44473
$ws.Cells.Item(16, 5).Value = 16
$ws.Cells.Item(16, 6).Value = "Fruta"
$ws.Cells.Item(16, 7).Value = 100104
$ws.Cells.Item(16, 8).Value = "Frutos de pepita"
$ws.Cells.Item(16, 9).Value = 100104005
$ws.Cells.Item(16, 10).Value = "Pera"
$ws.Cells.Item(16, 11).Value = "Packham's Triumph"
$ws.Cells.Item(16, 12).Value = "Primera"
$ws.Cells.Item(16, 13).Value = 60
$ws.Cells.Item(16, 14).Value = 9500
$ws.Cells.Item(16, 15).Value = 10000
$ws.Cells.Item(16, 16).Value = 9750
$ws.Cells.Item(16, 17).Value = "`$/caja 16 kilos empedrada"
$ws.Cells.Item(16, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(16, 19).Value = 609
$ws.Cells.Item(16, 20).Value = 16
